# 14th commit:
#  - Trim the "Login" sheet down to just Username/Password (drop the URL
#    column and the two duplicate data rows), keeping the surviving
#    mailto hyperlink.
#  - Add a new "DemoWebLogin" sheet (EmailID / Password) after "Login"
#    and make "Login" the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Login sheet: drop the extra rows (3 & 4) and the URL column (A),
#    shifting Username/Password into columns A/B.
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login")

# Remove the duplicate data rows first (rows 3 and 4).
$loginSheet.Range("A3:C4").Delete(-4162)

# Drop column A (URL), shifting Username/Password left.
$loginSheet.Range("A1:A2").Delete(-4159)

# The old hyperlinks collection (6 links across the old cells) no longer
# lines up with the new layout - clear it and re-add the single
# surviving link (B2, the Password/mailto one).
$loginSheet.Cells.Hyperlinks.Delete()
$loginSheet.Hyperlinks.Add($loginSheet.Range("B2"), "mailto:rmgy@9999")
# Re-applying the border keeps B2 on the worksheet's existing bordered
# "Hyperlink" style instead of a freshly minted one.
$loginSheet.Range("B2").Borders.LineStyle = 1

$loginSheet.Range("B6").Select()

# ---------------------------------------------------------------------
# 2. Add the new "DemoWebLogin" sheet after "Login".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$demoSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$demoSheet.Name = "DemoWebLogin"

$demoSheet.Range("A1").Value = "EmailID"
$demoSheet.Range("B1").Value = "Password"
$demoSheet.Range("A2").Value = "shekuemail@gmail.com"
$demoSheet.Range("B2").Value = '$Login123$'

$demoSheet.Columns.Item(1).ColumnWidth = 19.3
$demoSheet.Columns.Item(2).ColumnWidth = 10.65

# ---------------------------------------------------------------------
# 3. Make "Login" the active sheet/tab again.
# ---------------------------------------------------------------------
$loginSheet.Activate()
